$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shows")

# Add a new row (18) of data, mirroring the pattern of row 17 but with new show info
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "playlist-tue-12pm"
$ws.Range("C18").Value = "Playlist"
$ws.Range("D18").Value = "weekly"
$ws.Range("F18").Value = "No info available"
$ws.Range("I18").Value = "Unknown"
$ws.Range("J18").Value = "Unknown"
$ws.Range("O18").Value = "Unknown"
$ws.Range("R18").Value = "https://assets.podomatic.net/ts/28/84/58/info58064/640x640_17515723.jpg?1756449471"

# Match the wrap-text styling applied to column F in other rows
$ws.Range("F18").WrapText = $true

# Update sheet view: scroll/selection changes for mobile styling
$ws.Activate()
$ws.Range("A33").Select()
